$wb = $excel.ActiveWorkbook

# Row 20 on sheet ALC
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H20").Value = 2999
$ws.Range("I20").Value = 2999
$ws.Range("K20").Value = 2999
$ws.Range("M20").Value = -2769

# Row 35 on sheet ALC
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H35").Value = 2999
$ws.Range("I35").Value = 2999
$ws.Range("K35").Value = 2999
$ws.Range("M35").Value = -2620

# Row 55 on sheet ALC
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H55").Value = 133.05556
$ws.Range("I55").Value = 89.92308
$ws.Range("K55").Value = 89.92308
$ws.Range("M55").Value = 124.07692

# Row 74 on sheet ALC
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H74").Value = 5119.2666
$ws.Range("I74").Value = 5199.2144
$ws.Range("J74").Value = 4000
$ws.Range("K74").Value = 5199.2144
$ws.Range("L74").Value = 4000
$ws.Range("M74").Value = -4263.2144
$ws.Range("N74").Value = -5872

# Row 77 on sheet ALC
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H77").Value = 5119.2666
$ws.Range("I77").Value = 5199.2144
$ws.Range("J77").Value = 4000
$ws.Range("K77").Value = 25996.072
$ws.Range("L77").Value = 20000
$ws.Range("M77").Value = -21316.072
$ws.Range("N77").Value = -29360

# Row 92 on sheet ALC
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H92").Value = 41667340
$ws.Range("I92").Value = 55556304
$ws.Range("K92").Value = 55556304
$ws.Range("M92").Value = -55555056

# Row 100 on sheet ALC
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H100").Value = 1667.8096
$ws.Range("I100").Value = 1502.6
$ws.Range("J100").Value = 1719.4375
$ws.Range("K100").Value = 1502.6
$ws.Range("L100").Value = 1719.4375
$ws.Range("M100").Value = -961.5999999999999
$ws.Range("N100").Value = -2801.4375

# Row 132 on sheet ALC
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H132").Value = 4373.8335
$ws.Range("I132").Value = 2611.822
$ws.Range("J132").Value = 11940.117
$ws.Range("K132").Value = 7835.466
$ws.Range("L132").Value = 35820.351
$ws.Range("M132").Value = -5305.466
$ws.Range("N132").Value = -40880.351

# Row 135 on sheet ALC
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H135").Value = 909.8889
$ws.Range("I135").Value = 460.2258
$ws.Range("K135").Value = 4142.0322
$ws.Range("M135").Value = -1607.0322

# Row 137 on sheet ALC
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H137").Value = 3932.9167
$ws.Range("I137").Value = 5682.9165
$ws.Range("K137").Value = 17048.7495
$ws.Range("M137").Value = -14498.7495

# Row 138 on sheet ALC
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H138").Value = 11907225
$ws.Range("I138").Value = 30304978
$ws.Range("J138").Value = 2797.3726
$ws.Range("K138").Value = 90914934
$ws.Range("L138").Value = 8392.1178
$ws.Range("M138").Value = -90909794
$ws.Range("N138").Value = -18672.1178

# Row 141 on sheet ALC
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H141").Value = 12613.473
$ws.Range("I141").Value = 13926.259
$ws.Range("K141").Value = 41778.777
$ws.Range("M141").Value = -36598.777

# Row 32 on sheet ARM
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 4644.7607
$ws.Range("I32").Value = 4575.8555
$ws.Range("K32").Value = 4575.8555
$ws.Range("M32").Value = -4288.8555

# Row 102 on sheet ARM
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H102").Value = 3848.5
$ws.Range("I102").Value = 3903.4167
$ws.Range("J102").Value = 3738.6667
$ws.Range("K102").Value = 3903.4167
$ws.Range("L102").Value = 3738.6667
$ws.Range("M102").Value = -2281.4167
$ws.Range("N102").Value = -6982.6667

# Row 122 on sheet ARM
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H122").Value = 2585.1072
$ws.Range("I122").Value = 2064.04
$ws.Range("J122").Value = 6927.3335
$ws.Range("K122").Value = 6192.12
$ws.Range("L122").Value = 20782.0005
$ws.Range("M122").Value = -3742.12
$ws.Range("N122").Value = -25682.0005

# Row 26 on sheet BSM
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H26").Value = 31310
$ws.Range("I26").Value = 29465
$ws.Range("J26").Value = 35000
$ws.Range("K26").Value = 29465
$ws.Range("L26").Value = 35000
$ws.Range("M26").Value = -29173
$ws.Range("N26").Value = -35584

# Row 82 on sheet BSM
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H82").Value = 33959
$ws.Range("I82").Value = 26189.25
$ws.Range("K82").Value = 26189.25
$ws.Range("M82").Value = -25806.25

# Row 85 on sheet BSM
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H85").Value = 33959
$ws.Range("I85").Value = 26189.25
$ws.Range("K85").Value = 26189.25
$ws.Range("M85").Value = -24863.25

# Row 105 on sheet BSM
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H105").Value = 2914.6
$ws.Range("I105").Value = 2893.375
$ws.Range("J105").Value = 2999.5
$ws.Range("K105").Value = 2893.375
$ws.Range("L105").Value = 2999.5
$ws.Range("M105").Value = -1146.375
$ws.Range("N105").Value = -6493.5

# Row 134 on sheet BSM
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H134").Value = 3910.9875
$ws.Range("I134").Value = 2919.4126
$ws.Range("J134").Value = 7381.5
$ws.Range("K134").Value = 8758.237800000001
$ws.Range("L134").Value = 22144.5
$ws.Range("M134").Value = -6223.237800000001
$ws.Range("N134").Value = -27214.5

# Row 58 on sheet CRP
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H58").Value = 1130.6666
$ws.Range("I58").Value = 1007.5185
$ws.Range("K58").Value = 1007.5185
$ws.Range("M58").Value = -804.5185

# Row 69 on sheet CRP
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H69").Value = 22896.4
$ws.Range("I69").Value = 16373
$ws.Range("K69").Value = 16373
$ws.Range("M69").Value = -15624

# Row 72 on sheet CRP
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H72").Value = 22896.4
$ws.Range("I72").Value = 16373
$ws.Range("K72").Value = 49119
$ws.Range("M72").Value = -45375

# Row 105 on sheet CRP
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H105").Value = 4079
$ws.Range("I105").Value = 1747.2
$ws.Range("J105").Value = 11851.667
$ws.Range("K105").Value = 1747.2
$ws.Range("L105").Value = 11851.667
$ws.Range("M105").Value = -0.2000000000000455
$ws.Range("N105").Value = -15345.667

# Row 132 on sheet CRP
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H132").Value = 7654.7715
$ws.Range("I132").Value = 3045.5417
$ws.Range("K132").Value = 9136.625100000001
$ws.Range("M132").Value = -6606.625100000001

# Row 136 on sheet CRP
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H136").Value = 1130.6666
$ws.Range("I136").Value = 1007.5185
$ws.Range("K136").Value = 3022.5555
$ws.Range("M136").Value = -472.5554999999999

# Row 9 on sheet CUL
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H9").Value = 1102820.2
$ws.Range("I9").Value = 5500000.5
$ws.Range("J9").Value = 3525.125
$ws.Range("K9").Value = 16500001.5
$ws.Range("L9").Value = 10575.375
$ws.Range("M9").Value = -16499777.5
$ws.Range("N9").Value = -11023.375

# Row 26 on sheet CUL
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H26").Value = 140.3
$ws.Range("I26").Value = 100.5
$ws.Range("K26").Value = 301.5
$ws.Range("M26").Value = -13.5

# Row 46 on sheet CUL
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H46").Value = 273099650
$ws.Range("I46").Value = 299.66666
$ws.Range("J46").Value = 324305760
$ws.Range("K46").Value = 898.9999799999999
$ws.Range("L46").Value = 972917280
$ws.Range("M46").Value = -807.9999799999999
$ws.Range("N46").Value = -972917462

# Row 64 on sheet CUL
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H64").Value = 3106.6
$ws.Range("I64").Value = 1999
$ws.Range("J64").Value = 3185.7144
$ws.Range("K64").Value = 5997
$ws.Range("L64").Value = 9557.143199999999
$ws.Range("M64").Value = -5727
$ws.Range("N64").Value = -10097.1432

# Row 67 on sheet CUL
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H67").Value = 3106.6
$ws.Range("I67").Value = 1999
$ws.Range("J67").Value = 3185.7144
$ws.Range("K67").Value = 5997
$ws.Range("L67").Value = 9557.143199999999
$ws.Range("M67").Value = -5061
$ws.Range("N67").Value = -11429.1432

# Row 88 on sheet CUL
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H88").Value = 6262
$ws.Range("J88").Value = 8182.6665
$ws.Range("L88").Value = 24547.9995
$ws.Range("N88").Value = -25403.9995

# Row 91 on sheet CUL
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H91").Value = 6262
$ws.Range("J91").Value = 8182.6665
$ws.Range("L91").Value = 24547.9995
$ws.Range("N91").Value = -27511.9995

# Row 43 on sheet GSM
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H43").Value = 0
$ws.Range("I43").Value = 0
$ws.Range("K43").Value = 0
$ws.Range("M43").ClearContents()

# Row 57 on sheet GSM
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H57").Value = 30212.666
$ws.Range("I57").Value = 0
$ws.Range("K57").Value = 0
$ws.Range("M57").ClearContents()

# Row 122 on sheet GSM
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H122").Value = 6035.875
$ws.Range("I122").Value = 5473.5
$ws.Range("K122").Value = 16420.5
$ws.Range("M122").Value = -13970.5

# Row 134 on sheet GSM
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H134").Value = 304000
$ws.Range("J134").Value = 304000
$ws.Range("L134").Value = 912000
$ws.Range("N134").Value = -917070

# Row 135 on sheet GSM
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H135").Value = 135058.88
$ws.Range("J135").Value = 135058.88
$ws.Range("L135").Value = 135058.88
$ws.Range("N135").Value = -145198.88

# Row 58 on sheet LTW
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H58").Value = 48583
$ws.Range("I58").Value = 35374.5
$ws.Range("K58").Value = 35374.5
$ws.Range("M58").Value = -35114.5

# Row 70 on sheet LTW
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H70").Value = 25000
$ws.Range("J70").Value = 25000
$ws.Range("L70").Value = 25000
$ws.Range("N70").Value = -25540

# Row 73 on sheet LTW
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H73").Value = 25000
$ws.Range("J73").Value = 25000
$ws.Range("L73").Value = 25000
$ws.Range("N73").Value = -26872

# Row 105 on sheet LTW
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H105").Value = 105236
$ws.Range("J105").Value = 105236
$ws.Range("L105").Value = 105236
$ws.Range("N105").Value = -112224

# Row 51 on sheet WVR
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H51").Value = 10999
$ws.Range("I51").Value = 10999
$ws.Range("K51").Value = 10999
$ws.Range("M51").Value = -10489

# Row 52 on sheet WVR
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H52").Value = 13989.8
$ws.Range("I52").Value = 9987.5
$ws.Range("K52").Value = 9987.5
$ws.Range("M52").Value = -9761.5
